$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3

$ws.Range("C4").Select()
